# Add "Project Ideas" and "Registration Codes" content to the SAM TODO
# spreadsheet, matching the "add project ideas and registration codes to
# todo list spreadsheet" commit.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- rename the three sheets -------------------------------------------
$ws1.Name = "To Do"
$ws2.Name = "Project Ideas"
$ws3.Name = "Registration Codes"

# --- "Project Ideas" sheet ----------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 44.6
$ws2.Range("A1").Value = "Customizable reports"

# --- "Registration Codes" sheet -----------------------------------------
# Header row (order in which the original author typed the headers -
# column B's header was added after the first code row was started).
$ws3.Range("A1").Value = "Code"
$ws3.Range("C1").Value = "Version"
$ws3.Range("D1").Value = "Name"
$ws3.Range("E1").Value = "Email address"
$ws3.Range("F1").Value = "Country"

$ws3.Range("A2").Value = "CB18B612-F85E-47C9-AB54-AFBCD468BB4B"

$ws3.Range("B1").Value = "Date"

$ws3.Range("C2").NumberFormat = "@"
$ws3.Range("C2").Value = "2014.9.30"
$ws3.Range("C2").ClearFormats()

$ws3.Range("D2").Value = "Christina Schall"
$ws3.Range("E2").Value = "christina.schall@schottsolar.com"
$ws3.Range("F2").Value = "Germany"

# bold header row
$ws3.Range("A1:F1").Font.Bold = $true

# registration date (10/9/2014) as a real date value
$ws3.Range("B2").Value = Get-Date -Year 2014 -Month 10 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws3.Range("B2").NumberFormat = "m/d/yyyy"

# mailto hyperlink on the email address cell
$ws3.Hyperlinks.Add($ws3.Range("E2"), "mailto:christina.schall@schottsolar.com") | Out-Null

# column widths
$ws3.Columns.Item(1).ColumnWidth = 44.1
$ws3.Columns.Item(2).ColumnWidth = 12.71
$ws3.Columns.Item(3).ColumnWidth = 10.86
$ws3.Columns.Item(4).ColumnWidth = 17.71
$ws3.Columns.Item(5).ColumnWidth = 43.57
$ws3.Columns.Item(6).ColumnWidth = 21.86

# freeze the header row and leave the active cell on E15
$ws3.Activate()
$ws3.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws3.Range("E15").Select() | Out-Null

# restore the original active sheet / scroll position
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 22
